$wb = $excel.ActiveWorkbook

# Append a new sheet for 2025-07-11 at the end of the workbook (ranking data).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "2025-07-11"

# Header row
$newSheet.Cells.Item(1, 1).Value = "rank"
$newSheet.Cells.Item(1, 2).Value = "title"
$newSheet.Cells.Item(1, 3).Value = "author"
$newSheet.Cells.Item(1, 4).Value = "latest_episode"
$newSheet.Range("A1:D1").Font.Bold = $true

# Ranking rows (rank, title, author, latest_episode)
$newSheet.Cells.Item(2, 1).Value = 1
$newSheet.Cells.Item(2, 2).Value = '不徳のギルド'
$newSheet.Cells.Item(2, 3).Value = '河添太一'
$newSheet.Cells.Item(2, 4).Value = '第９６話：分福'

$newSheet.Cells.Item(3, 1).Value = 2
$newSheet.Cells.Item(3, 2).Value = '辺境モブ貴族のウチに嫁いできた悪役令嬢が、めちゃくちゃできる良い嫁なんだが？'
$newSheet.Cells.Item(3, 3).Value = 'tera(原作) 朝倉はやて(作画) 徹田(キャラクター原案)'
$newSheet.Cells.Item(3, 4).Value = '第9話'

$newSheet.Cells.Item(4, 1).Value = 3
$newSheet.Cells.Item(4, 2).Value = '落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～'
$newSheet.Cells.Item(4, 3).Value = '村上よしゆき 茨木野 あるてら'
$newSheet.Cells.Item(4, 4).Value = '第４０話　勇者、聖女と元聖騎士と再会し、魚人を追っ払う（４）'

$newSheet.Cells.Item(5, 1).Value = 4
$newSheet.Cells.Item(5, 2).Value = '魔王になったので、ダンジョン造って人外娘とほのぼのする'
$newSheet.Cells.Item(5, 3).Value = '遠野ノオト(作画) 流優(原作) だぶ竜(キャラクター原案)'
$newSheet.Cells.Item(5, 4).Value = '第1話後半'

$newSheet.Cells.Item(6, 1).Value = 5
$newSheet.Cells.Item(6, 2).Value = '俺以外誰も採取できない素材なのに「素材採取率が低い」とパワハラする幼馴染錬金術師と絶縁した専属魔導士、辺境の町でスローライフを送りたい。'
$newSheet.Cells.Item(6, 3).Value = '狐御前(原作) 西岡知三(作画) ＮＯＣＯ(キャラクター原案)'
$newSheet.Cells.Item(6, 4).Value = '第23話-2'

$newSheet.Cells.Item(7, 1).Value = 6
$newSheet.Cells.Item(7, 2).Value = 'バーサス'
$newSheet.Cells.Item(7, 3).Value = '原作：ONE 漫画：あずま京太郎 構成：bose'
$newSheet.Cells.Item(7, 4).Value = '第26話 惨事（1）'

$newSheet.Cells.Item(8, 1).Value = 7
$newSheet.Cells.Item(8, 2).Value = 'ゲーム悪役貴族に転生した俺は、チート筋肉で無双する'
$newSheet.Cells.Item(8, 3).Value = '昼行燈（原作） しいたけ元帥（漫画）'
$newSheet.Cells.Item(8, 4).Value = '第25話'

$newSheet.Cells.Item(9, 1).Value = 8
$newSheet.Cells.Item(9, 2).Value = '最弱貴族に転生したので悪役たちを集めてみた'
$newSheet.Cells.Item(9, 3).Value = '空野進 sorani ファルまろ'
$newSheet.Cells.Item(9, 4).Value = '第10話　最弱貴族、部下を信じる（１）'

$newSheet.Cells.Item(10, 1).Value = 9
$newSheet.Cells.Item(10, 2).Value = '俺は全てを【パリイ】する　～逆勘違いの世界最強は冒険者になりたい～'
$newSheet.Cells.Item(10, 3).Value = '原作：鍋敷・カワグチ 漫画：KRSG'
$newSheet.Cells.Item(10, 4).Value = '第24話'

$newSheet.Cells.Item(11, 1).Value = 10
$newSheet.Cells.Item(11, 2).Value = '無職の英雄　別にスキルなんか要らなかったんだが'
$newSheet.Cells.Item(11, 3).Value = '原作：九頭七尾・上田夢人 漫画：名苗秋緒'
$newSheet.Cells.Item(11, 4).Value = '第50話'

$newSheet.Cells.Item(12, 1).Value = 11
$newSheet.Cells.Item(12, 2).Value = 'ディーふらぐ！'
$newSheet.Cells.Item(12, 3).Value = '春野友矢(著者)'
$newSheet.Cells.Item(12, 4).Value = '第171話'

$newSheet.Cells.Item(13, 1).Value = 12
$newSheet.Cells.Item(13, 2).Value = '転生したらスライムだった件 異聞 ～魔国暮らしのトリニティ～'
$newSheet.Cells.Item(13, 3).Value = '伏瀬 戸野タエ みっつばー'
$newSheet.Cells.Item(13, 4).Value = '第106話　開国祭開幕［その１］'

$newSheet.Cells.Item(14, 1).Value = 13
$newSheet.Cells.Item(14, 2).Value = '双子まとめて『カノジョ』にしない？'
$newSheet.Cells.Item(14, 3).Value = '飴色みそ(漫画) 白井ムク(原作) 千種みのり(キャラクター原案)'
$newSheet.Cells.Item(14, 4).Value = '第13話①'

$newSheet.Cells.Item(15, 1).Value = 14
$newSheet.Cells.Item(15, 2).Value = 'クラスメイトの元アイドルが、とにかく挙動不審なんです。'
$newSheet.Cells.Item(15, 3).Value = 'となりける(作画) こりんさん(原作) ｋｒ木(キャラクター原案) マイクロマガジン社(監修)'
$newSheet.Cells.Item(15, 4).Value = '第24話-2'

$newSheet.Cells.Item(16, 1).Value = 15
$newSheet.Cells.Item(16, 2).Value = '35歳独身山田、異世界村に理想のセカンドハウスを作りたい　～異世界と現実のいいとこどりライフ～'
$newSheet.Cells.Item(16, 3).Value = '出雲大吉(原作) 西尾洋一(作画) ゆのひと(キャラクター原案)'
$newSheet.Cells.Item(16, 4).Value = '第1話④'

$newSheet.Cells.Item(17, 1).Value = 16
$newSheet.Cells.Item(17, 2).Value = '貧乏騎士に嫁入りしたはずが!? ~野人令嬢は皇太子妃になっても竜を狩りたい~'
$newSheet.Cells.Item(17, 3).Value = '漫画：夏川そぞろ 原作：宮前葵 キャラクター原案：ののまろ'
$newSheet.Cells.Item(17, 4).Value = '第10話③皇太子妃（仮）'

$newSheet.Cells.Item(18, 1).Value = 17
$newSheet.Cells.Item(18, 2).Value = '幼馴染のS級パーティーから追放された聖獣使い。万能支援魔法と仲間を増やして最強へ！'
$newSheet.Cells.Item(18, 3).Value = '黒田高祥(作画) かなりつ(原作) 転(キャラクター原案)'
$newSheet.Cells.Item(18, 4).Value = '第51話-1'

$newSheet.Cells.Item(19, 1).Value = 18
$newSheet.Cells.Item(19, 2).Value = 'となりの席のヤツがそういう目で見てくる'
$newSheet.Cells.Item(19, 3).Value = 'mmk'
$newSheet.Cells.Item(19, 4).Value = '第40話 誘惑'

$newSheet.Cells.Item(20, 1).Value = 19
$newSheet.Cells.Item(20, 2).Value = '転生してハイエルフになりましたが、スローライフは１２０年で飽きました'
$newSheet.Cells.Item(20, 3).Value = '原作：らる鳥・しあびす 漫画：成田コウ'
$newSheet.Cells.Item(20, 4).Value = '第40話'

$newSheet.Cells.Item(21, 1).Value = 20
$newSheet.Cells.Item(21, 2).Value = '【悲報】清楚系で売っていた底辺配信者、うっかり配信を切り忘れたままSS級モンスターを拳で殴り飛ばしてしまう'
$newSheet.Cells.Item(21, 3).Value = 'アトハ NEO草野 pupps'
$newSheet.Cells.Item(21, 4).Value = '第４話　【悲報】ご乱行⁉ ダンジョンシーカー・アカデミー！（３）'

$newSheet.Cells.Item(22, 1).Value = 21
$newSheet.Cells.Item(22, 2).Value = '実家に帰ったら甘やかされ生活が始まりました'
$newSheet.Cells.Item(22, 3).Value = '漫画：幹藻ねずみ 原作：月夜乃古狸 キャラクター原案：うなさか'
$newSheet.Cells.Item(22, 4).Value = '第23話後半'

$newSheet.Cells.Item(23, 1).Value = 22
$newSheet.Cells.Item(23, 2).Value = '戦隊レッド 異世界で冒険者になる'
$newSheet.Cells.Item(23, 3).Value = '中吉虎吉'
$newSheet.Cells.Item(23, 4).Value = '第41話「悠戯のブイダラ（前編）」'

$newSheet.Cells.Item(24, 1).Value = 23
$newSheet.Cells.Item(24, 2).Value = 'クラスメイトの美少女四人に頼まれたので、VRMMO内で専属料理人をはじめました'
$newSheet.Cells.Item(24, 3).Value = '斗樹稼多利(原作) 幾夜大黒堂(漫画) 中林ずん(キャラクター原案)'
$newSheet.Cells.Item(24, 4).Value = '第6話'

$newSheet.Cells.Item(25, 1).Value = 24
$newSheet.Cells.Item(25, 2).Value = 'ヘルモード ～やり込み好きのゲーマーは廃設定の異世界で無双する～ はじまりの召喚士'
$newSheet.Cells.Item(25, 3).Value = '原作：ハム男・藻 漫画：鉄田猿児'
$newSheet.Cells.Item(25, 4).Value = 'GAME 081　戦姫'

$newSheet.Cells.Item(26, 1).Value = 25
$newSheet.Cells.Item(26, 2).Value = 'フルメタル・パニック！　Family'
$newSheet.Cells.Item(26, 3).Value = '賀東招二(原作) 神反ヲ鬚(作画) 四季童子(キャラクター原案)'
$newSheet.Cells.Item(26, 4).Value = '第6話　東京都江東区のタワマン39階②-1'

$newSheet.Cells.Item(27, 1).Value = 26
$newSheet.Cells.Item(27, 2).Value = 'クラスで２番目に可愛い女の子と友だちになった'
$newSheet.Cells.Item(27, 3).Value = '尾野凛(漫画) たかた(原作) 日向あずり(キャラクター原案)'
$newSheet.Cells.Item(27, 4).Value = '第33話②'

$newSheet.Cells.Item(28, 1).Value = 27
$newSheet.Cells.Item(28, 2).Value = '規格外のダンジョン攻略者、実は異世界帰りの元勇者'
$newSheet.Cells.Item(28, 3).Value = '作画：やまざき君 原作：榊与一'
$newSheet.Cells.Item(28, 4).Value = '第4話(1)'

$newSheet.Cells.Item(29, 1).Value = 28
$newSheet.Cells.Item(29, 2).Value = 'アザミヤコを好きになる'
$newSheet.Cells.Item(29, 3).Value = 'ユニティコング(原作) ツノニガウ(作画)'
$newSheet.Cells.Item(29, 4).Value = '第8話'

$newSheet.Cells.Item(30, 1).Value = 29
$newSheet.Cells.Item(30, 2).Value = '俺の死亡フラグが留まるところを知らない'
$newSheet.Cells.Item(30, 3).Value = '漫画：乙須ミツヤ 原作：泉'
$newSheet.Cells.Item(30, 4).Value = 'フラグ68 フリエリ'

$newSheet.Cells.Item(31, 1).Value = 30
$newSheet.Cells.Item(31, 2).Value = '男女比1：5の世界でも普通に生きられると思った？　～激重感情な彼女たちが無自覚男子に翻弄されたら～'
$newSheet.Cells.Item(31, 3).Value = '三藤 孝太郎(原作) 桃季憂(漫画) jimmy(キャラクター原案)'
$newSheet.Cells.Item(31, 4).Value = '第9話-1'

$newSheet.Cells.Item(32, 1).Value = 31
$newSheet.Cells.Item(32, 2).Value = '今日も絵に描いた餅が美味い@COMIC'
$newSheet.Cells.Item(32, 3).Value = '漫画：梅渡飛鳥 原作：もちもち物質 キャラクター原案：転'
$newSheet.Cells.Item(32, 4).Value = '第41話'

$newSheet.Cells.Item(33, 1).Value = 32
$newSheet.Cells.Item(33, 2).Value = 'ギャルゲーマーに褒められたい'
$newSheet.Cells.Item(33, 3).Value = 'げしゅまろ(著者)'
$newSheet.Cells.Item(33, 4).Value = '40話'

$newSheet.Cells.Item(34, 1).Value = 33
$newSheet.Cells.Item(34, 2).Value = '斎藤義龍に生まれ変わったので、織田信長に国譲りして長生きするのを目指します！'
$newSheet.Cells.Item(34, 3).Value = '巽未頼 田村ゆうき マキムラシュンスケ'
$newSheet.Cells.Item(34, 4).Value = '第71話「日々の積み重ねこそ」'

$newSheet.Cells.Item(35, 1).Value = 34
$newSheet.Cells.Item(35, 2).Value = '実は俺、最強でした？'
$newSheet.Cells.Item(35, 3).Value = '原作：澄守 彩 漫画：高橋 愛'
$newSheet.Cells.Item(35, 4).Value = '第120話　四騎戦決勝戦!!・後編'

$newSheet.Cells.Item(36, 1).Value = 35
$newSheet.Cells.Item(36, 2).Value = '戦姫サバイバルサガ-異世界の運命をかけた無人島フジュン異性交遊-'
$newSheet.Cells.Item(36, 3).Value = 'OTOSAMA(著者)'
$newSheet.Cells.Item(36, 4).Value = '第17話'

$newSheet.Cells.Item(37, 1).Value = 36
$newSheet.Cells.Item(37, 2).Value = '無能なナナ'
$newSheet.Cells.Item(37, 3).Value = '原作 るーすぼーい 作画 古屋庵'
$newSheet.Cells.Item(37, 4).Value = '第78話 兄弟PART2'

$newSheet.Cells.Item(38, 1).Value = 37
$newSheet.Cells.Item(38, 2).Value = '殺されたらゾンビになったので、進化しまくって無双しようと思います'
$newSheet.Cells.Item(38, 3).Value = '漫画：朝ケ夜 原作：幸運ピエロ キャラクター原案：東西'
$newSheet.Cells.Item(38, 4).Value = '第15話(後半)暴走ドラゴンと魔剣②'

$newSheet.Cells.Item(39, 1).Value = 38
$newSheet.Cells.Item(39, 2).Value = 'ルパン三世 異世界の姫君（ネイバーワールドプリンセス）'
$newSheet.Cells.Item(39, 3).Value = 'モンキー・パンチ／エム・ピー・ワークス 内々けやき 佐伯庸介 白狼'
$newSheet.Cells.Item(39, 4).Value = '第100話：金毛羊の星空'

$newSheet.Cells.Item(40, 1).Value = 39
$newSheet.Cells.Item(40, 2).Value = 'モブ高生の俺でも冒険者になればリア充になれますか？'
$newSheet.Cells.Item(40, 3).Value = '原作：百均 漫画：さぎやまれん キャラクター原案：hai'
$newSheet.Cells.Item(40, 4).Value = '第29.5話'

$newSheet.Cells.Item(41, 1).Value = 40
$newSheet.Cells.Item(41, 2).Value = '異世界魔王と召喚少女の奴隷魔術'
$newSheet.Cells.Item(41, 3).Value = '原作：むらさきゆきや 漫画：福田直叶 キャラクター原案：鶴崎貴大'
$newSheet.Cells.Item(41, 4).Value = '第126話　戦争を終わらせてみるⅡ（中編）'

$newSheet.Cells.Item(42, 1).Value = 41
$newSheet.Cells.Item(42, 2).Value = '直径3cmの召喚陣<リミットリング>で「雑魚すら呼べない」と蔑まれた底辺召喚士が頂点に立つまで'
$newSheet.Cells.Item(42, 3).Value = '作画：まっつー 原作：空松蓮司'
$newSheet.Cells.Item(42, 4).Value = '第4話(1)'

$newSheet.Cells.Item(43, 1).Value = 42
$newSheet.Cells.Item(43, 2).Value = 'おねえさんと猫を飼う'
$newSheet.Cells.Item(43, 3).Value = '上杉響士郎(著者)'
$newSheet.Cells.Item(43, 4).Value = '第2話：おねえさんと猫の部屋'

$newSheet.Cells.Item(44, 1).Value = 43
$newSheet.Cells.Item(44, 2).Value = '六姫は神護衛に恋をする　～最強の守護騎士、転生して魔法学園に行く～'
$newSheet.Cells.Item(44, 3).Value = '漫画:加古山 寿 原案:朱月 十話 キャラクター原案:てつぶた'
$newSheet.Cells.Item(44, 4).Value = '第122話　声'

$newSheet.Cells.Item(45, 1).Value = 44
$newSheet.Cells.Item(45, 2).Value = '隣の席のヤンキー清水さんが髪を黒く染めてきた'
$newSheet.Cells.Item(45, 3).Value = '底花(原作) 真田若楓(漫画) ハム(キャラクター原案)'
$newSheet.Cells.Item(45, 4).Value = '第10話-1'

$newSheet.Cells.Item(46, 1).Value = 45
$newSheet.Cells.Item(46, 2).Value = '独身貴族は異世界を謳歌する ～結婚しない男の優雅なおひとりさまライフ～'
$newSheet.Cells.Item(46, 3).Value = '漫画：駒鳥 ひわ 原作：錬金王 キャラクター原案：三登 いつき'
$newSheet.Cells.Item(46, 4).Value = '第31話 独身貴族はヒラメが大事（4）'

$newSheet.Cells.Item(47, 1).Value = 46
$newSheet.Cells.Item(47, 2).Value = 'かつての暗殺者は来世で違う生き方をする'
$newSheet.Cells.Item(47, 3).Value = 'ツネ(漫画) 丘野優(原作) つなかわ(キャラクター原案)'
$newSheet.Cells.Item(47, 4).Value = '第4話①'

$newSheet.Cells.Item(48, 1).Value = 47
$newSheet.Cells.Item(48, 2).Value = '廃嫡王子の華麗なる逃亡劇 ~手段を選ばない最強クズ魔術師は自堕落に生きたい~'
$newSheet.Cells.Item(48, 3).Value = '出雲大吉(原作) 岡野むろ(作画) ゆのひと(キャラクター原案)'
$newSheet.Cells.Item(48, 4).Value = '第9話'

$newSheet.Cells.Item(49, 1).Value = 48
$newSheet.Cells.Item(49, 2).Value = 'スキル【万物支配】に目覚めたおっさんは、ダンジョンで生計を立てることにしました～無職から始める支配者無双～'
$newSheet.Cells.Item(49, 3).Value = '岸本和葉 原田 臙 シミズヒロノリ 吉武'
$newSheet.Cells.Item(49, 4).Value = '第4話　穴熊商店(前編)'

$newSheet.Cells.Item(50, 1).Value = 49
$newSheet.Cells.Item(50, 2).Value = '怠惰な悪辱貴族に転生した俺、シナリオをぶっ壊したら規格外の魔力で最凶になった'
$newSheet.Cells.Item(50, 3).Value = '菊池快晴(原作) 小田童馬(作画) 桑島黎音(キャラクター原案)'
$newSheet.Cells.Item(50, 4).Value = '重版決定記念イラスト'

$newSheet.Cells.Item(51, 1).Value = 50
$newSheet.Cells.Item(51, 2).Value = 'ゲーム内最強の『裏ボス』に転生したので、主人公の代わりに最速クリアを目指します！'
$newSheet.Cells.Item(51, 3).Value = '作画：こめぐ 原作：迅空也'
$newSheet.Cells.Item(51, 4).Value = '第4話(1)'
